$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty NAME cell (C115) with the company name,
# using a Text number format to match the other populated NAME cells
# in this column.
$ws.Range("C115").NumberFormat = "@"
$ws.Range("C115").Value = "ООО НКО ""Расчетные Решения»"

# Widen the NAME (C) and TYPE (D) columns so the newly-populated long
# text is readable; other columns keep their default width.
$ws.Columns.Item(3).ColumnWidth = 59.666666666666664
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666
